$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (price + 1h volume change) scraped by the GitHub Actions job.
# Numeric-looking text values are written with a leading apostrophe so Excel keeps them as text
# (matching the workbook's existing convention of storing Price/Volume columns as strings),
# exactly like typing them manually in the UI with a quote prefix.

$ws.Range('D2').Value = '30.120.35'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '1.906.68'
$ws.Range('E3').Value = '  -1.27%  '
$ws.Range('D4').Value = "'" + '0.9971'
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').Value = "'" + '0.7400'
$ws.Range('E5').Value = '  -2.12%  '
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('D7').Value = "'" + '0.9981'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').Value = "'" + '0.3110'
$ws.Range('E8').Value = '  -2.12%  '
$ws.Range('D9').Value = "'" + '26.50'
$ws.Range('E9').Value = '  -3.81%  '
$ws.Range('D10').Value = "'" + '0.06945'
$ws.Range('E10').Value = '  -0.80%  '
$ws.Range('D11').Value = "'" + '0.7757'
$ws.Range('E11').Value = '  -0.50%  '
$ws.Range('D12').Value = "'" + '0.07958'
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').Value = '1.895.44'
$ws.Range('E13').Value = '  -1.85%  '
$ws.Range('D14').Value = "'" + '5.251'
$ws.Range('E14').Value = '  -2.02%  '
$ws.Range('D15').Value = "'" + '91.93'
$ws.Range('E15').Value = '  -2.48%  '
$ws.Range('D16').Value = '30.103.58'
$ws.Range('E16').Value = '  -0.64%  '
$ws.Range('D17').Value = "'" + '14.19'
$ws.Range('E17').Value = '  -1.59%  '
$ws.Range('D18').Value = "'" + '5.829'
$ws.Range('E18').Value = '  +1.78%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = "'" + '241.04'
$ws.Range('E19').Value = '  -4.44%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = "'" + '0.000007808'
$ws.Range('E20').Value = '  -1.25%  '
$ws.Range('D21').Value = "'" + '0.9976'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('D22').Value = '2.149.11'
$ws.Range('E22').Value = '  -1.79%  '
$ws.Range('D23').Value = "'" + '0.9972'
$ws.Range('E23').Value = '  -0.25%  '
$ws.Range('D24').Value = "'" + '6.941'
$ws.Range('E24').Value = '  +4.03%  '
$ws.Range('D25').Value = "'" + '9.356'
$ws.Range('E25').Value = '  -1.38%  '
$ws.Range('D26').Value = "'" + '167.47'
$ws.Range('E26').Value = '  +0.97%  '
$ws.Range('D27').Value = "'" + '18.89'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').Value = "'" + '0.1281'
$ws.Range('E28').Value = '  -4.02%  '
$ws.Range('D29').Value = "'" + '2.054'
$ws.Range('E29').Value = '  -6.88%  '
$ws.Range('D30').Value = "'" + '1.350'
$ws.Range('E30').Value = '  -0.90%  '
$ws.Range('D31').Value = "'" + '1.540'
$ws.Range('E31').Value = '  +1.76%  '
$ws.Range('D32').Value = "'" + '4.316'
$ws.Range('E32').Value = '  -1.20%  '
$ws.Range('D33').Value = "'" + '4.080'
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('D34').Value = "'" + '0.05119'
$ws.Range('E34').Value = '  -0.73%  '
$ws.Range('D35').Value = "'" + '1.291'
$ws.Range('E35').Value = '  +0.60%  '
$ws.Range('D36').Value = "'" + '0.7390'
$ws.Range('E36').Value = '  -0.86%  '
$ws.Range('D37').Value = "'" + '2.709'
$ws.Range('E37').Value = '  -2.22%  '
$ws.Range('D38').Value = "'" + '0.01937'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').Value = "'" + '2.793'
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').Value = "'" + '6.315'
$ws.Range('E40').Value = '  -1.31%  '
$ws.Range('D41').Value = "'" + '74.72'
$ws.Range('E41').Value = '  -3.64%  '
$ws.Range('D42').Value = "'" + '0.4474'
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('D43').Value = "'" + '1.946'
$ws.Range('E43').Value = '  -0.88%  '
$ws.Range('D44').Value = "'" + '0.9993'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = "'" + '0.8360'
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('D46').Value = "'" + '7.775'
$ws.Range('E46').Value = '  +4.37%  '
$ws.Range('D47').Value = "'" + '101.28'
$ws.Range('E47').Value = '  +0.52%  '
$ws.Range('D48').Value = "'" + '9.858'
$ws.Range('E48').Value = '  +1.44%  '
$ws.Range('D49').Value = '2.054.53'
$ws.Range('D50').Value = "'" + '36.90'
$ws.Range('E50').Value = '  -1.11%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = "'" + '0.1185'
$ws.Range('E51').Value = '  +1.26%  '
